$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for numeric-looking price/volume strings so the
# COM layer does not coerce them (e.g. "235.00" -> 235, "0.0196" -> 1.96E-2),
# matching the workbook convention where these columns are stored as text.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '35.432.72'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').Value = '1.842.47'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '235.00'
$ws.Range('E5').Value = '  +4.40%  '
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = '43.83'
$ws.Range('E8').Value = '  +11.27%  '
$ws.Range('D9').Value = '0.312'
$ws.Range('E9').Value = '  +7.55%  '
$ws.Range('E10').Value = '  +3.62%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').Value = '1.842.39'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '11.31'
$ws.Range('E14').Value = '  +3.45%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.672'
$ws.Range('E15').Value = '  +5.84%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '4.75'
$ws.Range('E16').Value = '  +8.35%  '
$ws.Range('D17').Value = '35.449.41'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = '70.68'
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('D19').Value = '0.0₃0800'
$ws.Range('E19').Value = '  +4.00%  '
$ws.Range('D20').Value = '242.30'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = '11.99'
$ws.Range('E21').Value = '  +7.90%  '
$ws.Range('D22').Value = '4.65'
$ws.Range('E22').Value = '  +13.59%  '
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').Value = '171.01'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').Value = '7.91'
$ws.Range('E26').Value = '  +2.85%  '
$ws.Range('D27').Value = '17.73'
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = '1.61'
$ws.Range('E29').Value = '  +31.36%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = '3.343.07'
$ws.Range('E31').Value = '  +37.59%  '
$ws.Range('E32').Value = '  +9.34%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.97'
$ws.Range('E33').Value = '  +5.51%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '4.09'
$ws.Range('E34').Value = '  +6.08%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = '94.76'
$ws.Range('E36').Value = '  +14.86%  '
$ws.Range('D37').Value = '0.690'
$ws.Range('E37').Value = '  +7.50%  '
$ws.Range('D38').Value = '1.12'
$ws.Range('E38').Value = '  +6.65%  '
$ws.Range('D39').Value = '0.0196'
$ws.Range('E39').Value = '  +4.12%  '
$ws.Range('D40').Value = '15.37'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('D41').Value = '1.327.10'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('E42').Value = '  +7.01%  '
$ws.Range('D43').Value = '1.27'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = '6.29'
$ws.Range('E47').Value = '  +9.90%  '
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').Value = '101.72'
$ws.Range('E51').Value = '  -0.35%  '

# Restore default styling so the cells keep their original (unstyled) look
# now that the text values are committed.
$dataRange.Style = "Normal"

